# Updated: ut 22. 12. 2020
# Applies revised AgTests (H) / AgPosit (I) figures for existing rows 257-287
# and appends four new daily rows (288-291) for 2020-12-17 .. 2020-12-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised values for existing rows (columns H = AgTests, I = AgPosit) ---
$revisions = @{
    257 = @{ H = 5229;  I = 273  }
    258 = @{ H = 3860;  I = 218  }
    259 = @{ H = 6103;  I = 419  }
    260 = @{ H = 11479; I = 746  }
    261 = @{ H = 17378; I = 596  }
    262 = @{ H = 9215  }
    264 = @{ H = 41901; I = 845  }
    265 = @{ H = 17198 }
    266 = @{ H = 13085; I = 671  }
    267 = @{ H = 13894; I = 807  }
    268 = @{ H = 15137; I = 714  }
    269 = @{ H = 9451;  I = 411  }
    271 = @{ H = 42085; I = 1592 }
    272 = @{ H = 30464; I = 1641 }
    273 = @{ H = 26705; I = 1358 }
    274 = @{ H = 28086; I = 1309 }
    275 = @{ H = 28515; I = 1243 }
    276 = @{ H = 13858; I = 459  }
    277 = @{ H = 3252;  I = 125  }
    278 = @{ H = 29445; I = 2022 }
    279 = @{ H = 44046; I = 3121 }
    280 = @{ H = 35947; I = 2382 }
    281 = @{ H = 44974; I = 3255 }
    282 = @{ H = 46434; I = 2821 }
    283 = @{ H = 17661; I = 1043 }
    284 = @{ H = 1090;  I = 95   }
    285 = @{ H = 40066; I = 3374 }
    286 = @{ H = 54412; I = 4146 }
    287 = @{ H = 55695; I = 3766 }
}

foreach ($rowNum in $revisions.Keys) {
    $vals = $revisions[$rowNum]
    if ($vals.ContainsKey('H')) {
        $ws.Cells.Item($rowNum, 8).Value = $vals['H']
    }
    if ($vals.ContainsKey('I')) {
        $ws.Cells.Item($rowNum, 9).Value = $vals['I']
    }
}

# --- New rows appended at the bottom (A..I) ---
$newRows = @(
    @(288, 44182, 146124, 104560, 40124, 18022, 3991, 1440, 52854, 3878),
    @(289, 44183, 149275, 106361, 41404, 16197, 3151, 1510, 61403, 3500),
    @(290, 44184, 151336, 107828, 41953, 9821,  2061, 1555, 18289, 1500),
    @(291, 44185, 152555, 109807, 41130, 6291,  1219, 1618, 14351, 458)
)

foreach ($rowData in $newRows) {
    $rowNum = $rowData[0]
    $ws.Cells.Item($rowNum, 1).Value = $rowData[1]
    $ws.Cells.Item($rowNum, 2).Value = $rowData[2]
    $ws.Cells.Item($rowNum, 3).Value = $rowData[3]
    $ws.Cells.Item($rowNum, 4).Value = $rowData[4]
    $ws.Cells.Item($rowNum, 5).Value = $rowData[5]
    $ws.Cells.Item($rowNum, 6).Value = $rowData[6]
    $ws.Cells.Item($rowNum, 7).Value = $rowData[7]
    $ws.Cells.Item($rowNum, 8).Value = $rowData[8]
    $ws.Cells.Item($rowNum, 9).Value = $rowData[9]
}
